$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # "2021-Q4" - stays first

# Insert a new worksheet right after "2021-Q4" for the "2022-Q1" data
$wsNew = $wb.Worksheets.Add($null, $ws1)
$wsNew.Name = "2022-Q1"

# Seed header row + row layout/format by copying from the "2021-Q4" sheet,
# which already has the identical 8-column layout and styling. Column A of
# the header row is intentionally left blank (no A1 cell), same as on the
# "2021-Q4" sheet, so only B1:H1 is copied.
$ws1.Range("B1:H1").Copy($wsNew.Range("B1"))
for ($r = 2; $r -le 9; $r++) {
    $ws1.Range("A2:H2").Copy($wsNew.Range("A$r"))
}

function Set-TextValue($range, $text) {
    # Force the cell to keep its value as text (preserves leading zeros
    # and decimal-formatted numbers like "5.13" instead of silently
    # converting them to numeric cells), then drop back to the default
    # "Normal" style so no stray number-format styling is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2: fund 001075
$wsNew.Range("A2").Value = 0
Set-TextValue $wsNew.Range("B2") "001075"
$wsNew.Range("C2").Value = "宝盈转型动力灵活配置混合"
Set-TextValue $wsNew.Range("D2") "5.13"
Set-TextValue $wsNew.Range("E2") "86.64"
Set-TextValue $wsNew.Range("F2") "3.19"
Set-TextValue $wsNew.Range("G2") "0.1636"
$wsNew.Range("H2").Value = 9

# Row 3: fund 013714
$wsNew.Range("A3").Value = 1
Set-TextValue $wsNew.Range("B3") "013714"
$wsNew.Range("C3").Value = "方正富邦泰利12个月持有期混合A"
Set-TextValue $wsNew.Range("D3") "3.60"
Set-TextValue $wsNew.Range("E3") "20.66"
Set-TextValue $wsNew.Range("F3") "2.03"
Set-TextValue $wsNew.Range("G3") "0.0731"
$wsNew.Range("H3").Value = 3

# Row 4: fund 011501
$wsNew.Range("A4").Value = 2
Set-TextValue $wsNew.Range("B4") "011501"
$wsNew.Range("C4").Value = "方正富邦汇福一年定期开放灵活配置混合A"
Set-TextValue $wsNew.Range("D4") "4.76"
Set-TextValue $wsNew.Range("E4") "41.05"
Set-TextValue $wsNew.Range("F4") "1.51"
Set-TextValue $wsNew.Range("G4") "0.0719"
$wsNew.Range("H4").Value = 9

# Row 5: fund 008602
$wsNew.Range("A5").Value = 3
Set-TextValue $wsNew.Range("B5") "008602"
$wsNew.Range("C5").Value = "方正富邦新兴成长混合A"
Set-TextValue $wsNew.Range("D5") "1.85"
Set-TextValue $wsNew.Range("E5") "78.95"
Set-TextValue $wsNew.Range("F5") "2.30"
Set-TextValue $wsNew.Range("G5") "0.0426"
$wsNew.Range("H5").Value = 9

# Row 6: fund 582003
$wsNew.Range("A6").Value = 4
Set-TextValue $wsNew.Range("B6") "582003"
$wsNew.Range("C6").Value = "东吴配置优化灵活配置混合"
Set-TextValue $wsNew.Range("D6") "1.04"
Set-TextValue $wsNew.Range("E6") "90.74"
Set-TextValue $wsNew.Range("F6") "3.50"
Set-TextValue $wsNew.Range("G6") "0.0364"
$wsNew.Range("H6").Value = 6

# Row 7: fund 013715
$wsNew.Range("A7").Value = 5
Set-TextValue $wsNew.Range("B7") "013715"
$wsNew.Range("C7").Value = "方正富邦泰利12个月持有期混合C"
Set-TextValue $wsNew.Range("D7") "0.10"
Set-TextValue $wsNew.Range("E7") "20.66"
Set-TextValue $wsNew.Range("F7") "2.03"
Set-TextValue $wsNew.Range("G7") "0.0020"
$wsNew.Range("H7").Value = 3

# Row 8: fund 011502
$wsNew.Range("A8").Value = 6
Set-TextValue $wsNew.Range("B8") "011502"
$wsNew.Range("C8").Value = "方正富邦汇福一年定期开放灵活配置混合C"
Set-TextValue $wsNew.Range("D8") "0.09"
Set-TextValue $wsNew.Range("E8") "41.05"
Set-TextValue $wsNew.Range("F8") "1.51"
Set-TextValue $wsNew.Range("G8") "0.0014"
$wsNew.Range("H8").Value = 9

# Row 9: fund 008603
$wsNew.Range("A9").Value = 7
Set-TextValue $wsNew.Range("B9") "008603"
$wsNew.Range("C9").Value = "方正富邦新兴成长混合C"
Set-TextValue $wsNew.Range("D9") "0.04"
Set-TextValue $wsNew.Range("E9") "78.95"
Set-TextValue $wsNew.Range("F9") "2.30"
Set-TextValue $wsNew.Range("G9") "0.0009"
$wsNew.Range("H9").Value = 9

# --- Update the "总计" (totals) sheet: keep existing "2021-Q4" summary row,
# but move it down to row 3 and insert a new "2022-Q1" summary row above it.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 8
$wsTotal.Range("D2").Value = 0.39

$wsTotal.Range("A3").Value = 1
